$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: new time-recording entry ("10월 30일" / 14:00-17:00 / 0 interrupt / 180 delta / "Nodejs 강의 수강")
$ws.Range("A11").Value = "10월 30일"
$ws.Range("A11").Characters(4, 4).Font.Name = "돋움"
$ws.Range("A11").Characters(4, 4).Font.Size = 10
$ws.Range("A11").Characters(4, 4).Font.ColorIndex = -4105

$ws.Range("B11").Value = 0.58333333333333337
$ws.Range("C11").Value = 0.70833333333333337
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 180

$ws.Range("F11").Value = "Nodejs 강의 수강"
$ws.Range("F11").Characters(10, 3).Font.Name = "돋움"
$ws.Range("F11").Characters(10, 3).Font.Size = 10
$ws.Range("F11").Characters(10, 3).Font.ColorIndex = -4105

# Active selection moves to F12
$ws.Range("F12").Select() | Out-Null
